$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.177.73"
$ws.Range("E2").Value = "  -1.58%  "

$ws.Range("D3").Value = "2.270.35"
$ws.Range("E3").Value = "  -2.46%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.36"
$ws.Range("E5").Value = "  -2.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.07"
$ws.Range("E6").Value = "  -5.18%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.494"
$ws.Range("E7").Value = "  -2.86%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -2.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.37"
$ws.Range("E10").Value = "  -3.29%  "

$ws.Range("E11").Value = "  -0.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.21"
$ws.Range("E12").Value = "  -7.53%  "

$ws.Range("E13").Value = "  -0.13%  "

$ws.Range("E14").Value = "  -1.43%  "

$ws.Range("D15").Value = "2.621.45"
$ws.Range("E15").Value = "  -2.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.50"
$ws.Range("E16").Value = "  -1.86%  "

$ws.Range("D17").Value = "2.340.77"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.782"
$ws.Range("E18").Value = "  -5.40%  "

$ws.Range("D19").Value = "42.097.97"
$ws.Range("E19").Value = "  -1.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.73"
$ws.Range("E20").Value = "  +1.09%  "

$ws.Range("E21").Value = "  -1.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.99"
$ws.Range("E22").Value = "  -2.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.60"
$ws.Range("E23").Value = "  -3.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.45"
$ws.Range("E24").Value = "  -0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.96"
$ws.Range("E25").Value = "  -1.07%  "

$ws.Range("E26").Value = "  +0.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.44"
$ws.Range("E27").Value = "  -3.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.96"
$ws.Range("E28").Value = "  -5.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  +2.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.62"
$ws.Range("E30").Value = "  +4.67%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.00"
$ws.Range("E31").Value = "  -2.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.10"
$ws.Range("E32").Value = "  -1.36%  "

$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.90"
$ws.Range("E34").Value = "  -3.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.51"
$ws.Range("E35").Value = "  -2.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.54"
$ws.Range("E36").Value = "  -3.40%  "

$ws.Range("E37").Value = "  -4.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0686"
$ws.Range("E38").Value = "  -4.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.78"
$ws.Range("E39").Value = "  -3.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0986"
$ws.Range("E40").Value = "  -2.52%  "

$ws.Range("E41").Value = "  -2.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.72"
$ws.Range("E42").Value = "  -5.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.45"
$ws.Range("E43").Value = "  -4.60%  "

$ws.Range("D44").Value = "1.954.76"
$ws.Range("E44").Value = "  -2.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0277"
$ws.Range("E45").Value = "  -1.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.37"
$ws.Range("E46").Value = "  -6.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.52"
$ws.Range("E47").Value = "  -6.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("E48").Value = "  -3.37%  "

$ws.Range("D49").Value = "2.493.90"
$ws.Range("E49").Value = "  -2.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.08"
$ws.Range("E50").Value = "  -6.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.51"
$ws.Range("E51").Value = "  -3.52%  "
